$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Rename the "View" field header to "Cache"
$ws.Range("F1").Value = "Cache"

# Default value of the Cache field column should now be FALSE instead of TRUE
$ws.Range("F2:F15").Value = $false
